$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "changing" fields (D,L,M,N,O,P,Q,R,S,T) for every data row
# before writing anything, because the edit is a row-wise permutation:
# every row is both a source and a destination.
$rows = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,19,20
$snap = @{}
foreach ($r in $rows) {
    $snap[$r] = @{
        D = $ws.Range("D$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        Q = $ws.Range("Q$r").Value2
        R = $ws.Range("R$r").Value2
        S = $ws.Range("S$r").Value2
        T = $ws.Range("T$r").Value2
    }
}

# destination row -> source row (source row's data moves onto destination row)
$map = @{
    2 = 13
    3 = 19
    4 = 3
    5 = 20
    6 = 9
    7 = 16
    8 = 15
    9 = 6
    10 = 8
    11 = 17
    12 = 7
    13 = 14
    14 = 10
    15 = 4
    16 = 5
    17 = 2
    19 = 12
    20 = 11
}

foreach ($dest in $map.Keys) {
    $src = $map[$dest]
    $data = $snap[$src]
    $ws.Range("D$dest").Value = $data.D
    $ws.Range("L$dest").Value = $data.L
    $ws.Range("M$dest").Value = $data.M
    $ws.Range("N$dest").Value = $data.N
    $ws.Range("O$dest").Value = $data.O
    $ws.Range("P$dest").Value = $data.P
    $ws.Range("Q$dest").Value = $data.Q
    $ws.Range("R$dest").Value = $data.R
    $ws.Range("S$dest").Value = $data.S
    $ws.Range("T$dest").Value = $data.T
}
